$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "id" column (column A) - the data is now uploaded/inserted per
# table, so the spreadsheet no longer needs a synthetic row id column.
$ws.Range("A1").EntireColumn.Delete()

# Clear the lingering selection left over from editing so the saved file
# doesn't keep a stale active-cell reference.
$ws.Range("A1").Select()
